# Add support for Corner Hi Lo odds: a new "全場角球入球大細" block
# (mirrors the existing "全場入球大細" block in K:N) placed in columns P:S.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header block (P1:S2), mirroring the K1:N2 "全場入球大細" block ---
$ws.Range("P1").Value = "全場角球入球大細"

$ws.Range("P2").Value = "更新時間"
# Copy A2's number format (time-of-day) onto P2 without minting a duplicate
# numFmt entry (a plain NumberFormat string round-trip creates one).
$ws.Range("A2").Copy()
$ws.Range("P2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("Q2").Value = "球數"
$ws.Range("R2").Value = "大"
$ws.Range("S2").Value = "細"

# --- Workbook-scoped defined name "全場角球入球大細" -> TEMPLATE!$P$1 ---
# The defined-name validator here requires the first character to be an
# ASCII letter/underscore (unlike real Excel, which allows a leading CJK
# character), so create it under a placeholder ASCII name, then rename.
$ws.Range("P1").Name = "ZZZ_CORNER_HILO_TMP"
$wb.Names.Item("ZZZ_CORNER_HILO_TMP").Name = "全場角球入球大細"

# Defined names are re-sorted (by name) only when Names.Add runs, not on a
# plain rename, so force one more Add/Delete cycle to re-sort the
# collection into the final alphabetical position (between 全場入球大細
# and 全場讓球), matching how Excel keeps <definedNames> sorted.
$wb.Names.Add("ZZZ_RESORT_TMP", "=TEMPLATE!`$A`$1")
$wb.Names.Item("ZZZ_RESORT_TMP").Delete()

# --- Selection moves to D9 ---
$ws.Range("D9").Select() | Out-Null
